# Add SmartGit to demo slides
# Slide 9 (homework slide) second shape contains the bulleted list that
# starts with "Инсталирайте Git Extensions". The commit appends ", SmartGit"
# to that single run's text without altering any other formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Target only the exact characters of the existing run (27 characters,
# excluding the paragraph mark) so the replacement stays a single run with
# its original formatting intact.
$target = $tr.Characters(1, 27)
$target.Text = "Инсталирайте Git Extensions, SmartGit"
